# Latest push for 04Oct
# - Flip the "InGeo"/ActivateIncomingRouting-style "No" flag in column A
#   (row 2 on every sheet, plus row 3 on DDI3/DDI2) over to "Yes".
# - Move the remembered cell selection on each sheet.

$wb = $excel.ActiveWorkbook

$wsDDI  = $wb.Worksheets.Item("DDI")
$wsDDI3 = $wb.Worksheets.Item("DDI3")
$wsDDI2 = $wb.Worksheets.Item("DDI2")

# Sheet "DDI": A2 No -> Yes
$wsDDI.Range("A2").Value = "Yes"

# Sheet "DDI3": A2/A3 No -> Yes
$wsDDI3.Range("A2").Value = "Yes"
$wsDDI3.Range("A3").Value = "Yes"

# Sheet "DDI2": A2/A3 No -> Yes
$wsDDI2.Range("A2").Value = "Yes"
$wsDDI2.Range("A3").Value = "Yes"

# Restore each sheet's remembered selection (cursor position).
$wsDDI.Activate()
$wsDDI.Range("C10").Select() | Out-Null

$wsDDI3.Activate()
$wsDDI3.Range("D12").Select() | Out-Null

$wsDDI2.Activate()
$wsDDI2.Range("D21").Select() | Out-Null

# "DDI" was the originally active/selected tab - leave it active again.
$wsDDI.Activate()
